# [EXTRA SCRAPE] full data scraped for extra batting and bowling fields
#
# 22 new player-match rows are inserted at the top of the "ODI Batting
# Extra" sheet's data (ahead of MATCH_CODE 4524...), pushing the existing
# rows 2-21 down to rows 24-43. The sheet's used range grows from
# A1:F21 to A1:F43.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("ODI Batting Extra")

# Push the existing data rows (2-21) down by 22 rows, preserving their
# content/types exactly, and opening up rows 2-23 for the new records.
$ws.Rows("2:23").Insert()

# Inserting above row 24 pulls in neighbouring (header) formatting; reset
# the freshly opened rows back to the plain/default style used by the
# rest of the data rows before filling them in.
$ws.Range("A2:F23").ClearFormats()

# Columns: MATCH_CODE(A, text), BATTING_POSITION(B, number-or-blank),
# NUM_4(C, text-or-blank), NUM_6(D, text-or-blank),
# PERCENT_RUNS_OF_TOTAL(E, text-or-blank), MAN_OF_MATCH(F, text)
$rows = @(
    @(2, "4096", "3", "1", "0", "8.04%", "NO"),
    @(3, "4098", "", "", "", "", "NO"),
    @(4, "4099", "3", "8", "1", "29.68%", "NO"),
    @(5, "4130", "5", "2", "0", "6.23%", "NO"),
    @(6, "4133", "5", "2", "0", "10.95%", "NO"),
    @(7, "4135", "5", "", "", "", "NO"),
    @(8, "4359", "", "", "", "", "NO"),
    @(9, "4360", "5", "5", "0", "25.45%", "NO"),
    @(10, "4362", "", "", "", "", "NO"),
    @(11, "4385", "", "", "", "", "NO"),
    @(12, "4387", "4", "3", "4", "13.70%", "NO"),
    @(13, "4388", "4", "1", "0", "2.22%", "NO"),
    @(14, "4398", "5", "0", "0", "1.57%", "NO"),
    @(15, "4399", "4", "1", "0", "2.06%", "NO"),
    @(16, "4400", "4", "6", "1", "15.22%", "NO"),
    @(17, "4402", "", "", "", "", "NO"),
    @(18, "4406", "", "", "", "", "NO"),
    @(19, "4410", "4", "9", "0", "20.95%", "NO"),
    @(20, "4435", "", "", "", "", "NO"),
    @(21, "4436", "", "", "", "", "NO"),
    @(22, "4437", "4", "2", "0", "6.29%", "NO"),
    @(23, "4454", "4", "1", "0", "1.89%", "NO")
)

foreach ($row in $rows) {
    $r = $row[0]
    $matchCode = $row[1]
    $battingPos = $row[2]
    $num4 = $row[3]
    $num6 = $row[4]
    $pctRuns = $row[5]
    $manOfMatch = $row[6]

    # A: MATCH_CODE, always stored as text
    $ws.Cells.Item($r, 1).Value = "'" + $matchCode

    # B: BATTING_POSITION, numeric when present, blank (empty text) otherwise
    if ($battingPos -eq "") {
        $ws.Cells.Item($r, 2).Value = "'"
    } else {
        $ws.Cells.Item($r, 2).Value = [double]$battingPos
    }

    # C: NUM_4, stored as text when present, blank (empty text) otherwise
    if ($num4 -eq "") {
        $ws.Cells.Item($r, 3).Value = "'"
    } else {
        $ws.Cells.Item($r, 3).Value = "'" + $num4
    }

    # D: NUM_6, stored as text when present, blank (empty text) otherwise
    if ($num6 -eq "") {
        $ws.Cells.Item($r, 4).Value = "'"
    } else {
        $ws.Cells.Item($r, 4).Value = "'" + $num6
    }

    # E: PERCENT_RUNS_OF_TOTAL, stored as text when present, blank (empty text) otherwise
    if ($pctRuns -eq "") {
        $ws.Cells.Item($r, 5).Value = "'"
    } else {
        $ws.Cells.Item($r, 5).Value = "'" + $pctRuns
    }

    # F: MAN_OF_MATCH, always text (YES/NO)
    $ws.Cells.Item($r, 6).Value = "'" + $manOfMatch
}

# The leading apostrophes above force text storage but also leave a
# "quote prefix" style flag behind; strip it so the new rows end up with
# the same plain (unstyled) look as the rest of the data rows.
$ws.Range("A2:F23").ClearFormats()
